$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F10").Value = 5787
$ws.Range("F13").Value = 1440
$ws.Range("F14").Value = 618
$ws.Range("F16").Value = 398
$ws.Range("F19").Value = 4805
$ws.Range("F22").Value = 2426
$ws.Range("F35").Value = 33
$ws.Range("F37").Value = 1413
$ws.Range("F39").Value = 101
$ws.Range("F40").Value = 537
$ws.Range("F41").Value = 202
$ws.Range("F42").Value = 1671
$ws.Range("F45").Value = 91
$ws.Range("F48").Value = 41
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 11
$ws.Range("F8").Value = 394
$ws.Range("F9").Value = 288
$ws.Range("F11").Value = 66
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1679
$ws.Range("F7").Value = 553
$ws.Range("F10").Value = 1777
$ws.Range("F11").Value = 2304
$ws.Range("F12").Value = 737
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1679
$ws.Range("F6").Value = 553
$ws.Range("F11").Value = 2304
$ws.Range("F12").Value = 5787
$ws.Range("F13").Value = 737
$ws.Range("F15").Value = 1440
$ws.Range("F18").Value = 11
$ws.Range("F20").Value = 4805
$ws.Range("F21").Value = 2426
$ws.Range("F27").Value = 288
$ws.Range("F29").Value = 66
$ws.Range("F37").Value = 1413
$ws.Range("F39").Value = 537
$ws.Range("F42").Value = 202
$ws.Range("F44").Value = 1671
$ws.Range("F46").Value = 91
$ws.Range("F49").Value = 41
